$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.052.40"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.053.92"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.46"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.83"
$ws.Range("E7").Value = "  +7.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.31"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.387"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  +6.11%  "
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.38"
$ws.Range("E13").Value = "  +8.41%  "
$ws.Range("D14").Value = "2.358.64"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.806"
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("E16").Value = "  +6.76%  "
$ws.Range("D17").Value = "2.093.01"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "37.014.95"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.56"
$ws.Range("E19").Value = "  +14.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.64"
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("D21").Value = "0.0₃0917"
$ws.Range("E21").Value = "  +7.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.46"
$ws.Range("E22").Value = "  +3.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.17"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  +12.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.19"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.27"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.27"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("E31").Value = "  +4.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.75"
$ws.Range("E32").Value = "  +4.61%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.49"
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0889"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.116"
$ws.Range("E38").Value = "  +18.65%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.76"
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.79"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.60"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.53"
$ws.Range("E46").Value = "  +11.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  +5.37%  "
$ws.Range("D48").Value = "1.292.39"
$ws.Range("E48").Value = "  -3.02%  "
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.90"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.246.54"
$ws.Range("E51").Value = "  -1.76%  "
